# The workbook ships with a sheet literally named "cotizacion " (trailing
# space). Rename it to the clean "cotizacion" and make sure the sheet's
# Print Area (which Excel stores as the workbook-level defined name
# `_xlnm.Print_Area`) is re-pointed at the new, space-free sheet name so it
# no longer needs to be quoted.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("cotizacion ")
$ws.Name = "cotizacion"

# Re-assert the print area so Excel regenerates the defined name using the
# sheet's current (trimmed) name instead of the stale quoted reference.
$ws.PageSetup.PrintArea = '$A$1:$I$42'
